$p = $ppt.ActivePresentation

# --- Slide 3: "Requirements and Scope" content placeholder -----------------
# Remove the leading "Example - " prefix from the two RFC bullet lines.
$s3 = $p.Slides.Item(3)
$contentShape = $s3.Shapes.Item(2)
$contentTr = $contentShape.TextFrame.TextRange

$para12 = $contentTr.Paragraphs(12, 1)
$para12.Runs(1, 1).Text = "RFC 5357 (TWAMP Light) defined probe messages"

$para13 = $contentTr.Paragraphs(13, 1)
$para13.Runs(1, 1).Text = "RFC 8762 (STAMP) defined probe messages"

# --- Slide 11: "Example Provisioning Model" diagram rectangle --------------
$s11 = $p.Slides.Item(11)
$diagramShape = $s11.Shapes.Item(4)
$diagramTr = $diagramShape.TextFrame.TextRange

# Paragraph 11 originally reads:
#   "                            /                  \"
# It becomes three runs: "     Source/" + "Dest" + " UDP Ports  /                  \"
$para11 = $diagramTr.Paragraphs(11, 1)
$run1 = $para11.Runs(1, 1)
$run1.Text = "     Source/"
$run2 = $run1.InsertAfter("Dest")
$run3 = $run2.InsertAfter(" UDP Ports  /                  \")

# Paragraph 20, the figure caption, gains two extra leading spaces.
$para20 = $diagramTr.Paragraphs(20, 1)
$para20.Runs(1, 1).Text = "                     Figure 2: Example Provisioning Model"

# Editing the text re-triggers this autofit textbox's layout pass, which
# would otherwise grow the shape's stored height; the original deck's box
# size is unaffected by this edit, so restore it explicitly (4184415 EMU).
$diagramShape.Height = 4184415 / 12700.0
